$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("menu")

$rows = 4..13
$urls = @(
  "https://static.order.kfcku.com/images/items/lg/5587.jpg?v=gqGmqL",
  "https://brand-uk.assets.kfc.co.uk/2023-10/KFC4317~39047_W8_23_ZINGER_BURGER_MEAL_1200x800.jpg?VersionId=APSJn8Jb5_jlGX1O3FEjwHadslkrlfXa",
  "https://www.shutterstock.com/image-photo/brno-czechia-06162023-starbucks-coffee-600nw-2320284399.jpg",
  "https://www.starbucks.co.th/stb-media/2020/08/1.Espresso-1080.png",
  "https://s7d1.scene7.com/is/image/mcdonalds/Header_BigMac_832x472:1-3-product-tile-desktop?wid=763&hei=472&dpr=off",
  "https://www.sidechef.com/recipe/b4496acb-32b7-4b27-8d91-b6e7a3c21218.jpg?d=1408x1120",
  "https://www.keyingredient.com/uploads/media/be/80/8316c8fb81a8cb68a15047c2748cf376d525.jpg?1560436606",
  "https://www.shape.com/thmb/DJKMBoJgTF_lfPnzAI4cnLhGD5w=/1500x0/filters:no_upscale():max_bytes(150000):strip_icc()/hot-coffee-8079cd4654e84c07924b12d4f78b2cb5.jpg",
  "https://static.wikia.nocookie.net/burgerking/images/3/3c/Whopper_JR.jpg/revision/latest?cb=20190619155429",
  "https://i0.wp.com/www.theimpulsivebuy.com/wordpress/wp-content/uploads/2023/10/bkghopepchickfries.jpg?fit=1200%2C1200&ssl=1"
)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $row = $rows[$i]
  $url = $urls[$i]
  $cell = $ws.Range("E$row")
  $cell.Value = $url
  $ws.Hyperlinks.Add($cell, $url)
  $cell.Style = "Hyperlink"
}
